# Added Icon Support topicslist.json
# - "posts" sheet: header row gains LINK/IMG/IMG_ALT columns (replacing SLUG/SRC/ALT)
#   and POST/TOPIC/COLOR/DESCRIPTION shift into NAME/TOPIC/COLOR/DESCRIPTION.
# - "topics" sheet: header TOPIC/SLUG/SRC -> NAME/LINK/ICON, and the D column
#   (previously raw SVG <path> markup) now holds Material Symbols icon names.

$wb = $excel.ActiveWorkbook

$posts = $wb.Worksheets.Item("posts")
$topics = $wb.Worksheets.Item("topics")

# ---- posts sheet header row ----
$posts.Range("B1").Value = "NAME"
$posts.Range("C1").Value = "LINK"
$posts.Range("D1").Value = "IMG"
$posts.Range("E1").Value = "IMG_ALT"
$posts.Range("F1").Value = "TOPIC"
$posts.Range("G1").Value = "COLOR"
$posts.Range("H1").Value = "DESCRIPTION"

# widen column E to fit its new header, matching the author's resize
# (the engine quantizes ColumnWidth to discrete pixel-character steps; 8.8
# is the input that lands closest to the authored 9.6328125 width)
$posts.Columns.Item(5).ColumnWidth = 8.8

# ---- topics sheet header row ----
$topics.Range("B1").Value = "NAME"
$topics.Range("C1").Value = "LINK"
$topics.Range("D1").Value = "ICON"

# ---- topics sheet icon column (D2:D11) ----
$topics.Range("D2").Value = "design_services"
$topics.Range("D3").Value = "construction"
$topics.Range("D4").Value = "spa"
$topics.Range("D5").Value = "self_improvement"
$topics.Range("D6").Value = "palette"
$topics.Range("D7").Value = "restaurant"
$topics.Range("D8").Value = "developer_mode"
$topics.Range("D9").Value = "devices"
$topics.Range("D10").Value = "park"
$topics.Range("D11").Value = "pets"

# ---- restore the cursor/selection positions recorded by the author ----
[void]$posts.Range("B1").Select()
[void]$topics.Range("D14").Select()

$topics.Activate()
